$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regression")

# New "Contact Us" / "List of Products" test case columns (D:F) for TC107 rows 19-20.
# Entered column-by-column (header then value) so new shared-string indices line up
# with the source workbook: Subject, Just a sample check, Message,
# Typing the Message…., FileName, C:\Users\Dileep_K\Desktop
$ws.Range("D19").Value = "Subject"
$ws.Range("D20").Value = "Just a sample check"

$ws.Range("E19").Value = "Message"
$ws.Range("E20").Value = "Typing the Message…."

$ws.Range("F19").Value = "FileName"
$ws.Range("F20").Value = "C:\Users\Dileep_K\Desktop"

# Widen columns E (5) and F (6) so the new content fits/best-fits.
$ws.Columns.Item(5).ColumnWidth = 34
$ws.Columns.Item(6).ColumnWidth = 42

# Move the active selection to the newly filled-in cell F20.
$ws.Range("F20").Select()
